# Change the table style on the "Sources of finance" table (slide 6) from the
# default "Table_0" style ({DCC4BDE6-7913-4E3F-87A6-BE83CE810A61}) to the
# built-in PowerPoint table style {F5A0190C-B820-4218-89FA-BDFADE48523E}.

$p = $ppt.ActivePresentation
$targetStyleId = "{F5A0190C-B820-4218-89FA-BDFADE48523E}"

$s = $p.Slides.Item(6)

$tableShape = $null
for ($j = 1; $j -le $s.Shapes.Count; $j++) {
    $candidate = $s.Shapes.Item($j)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}

if ($tableShape -eq $null) {
    $tableShape = $s.Shapes.Item(2)
}

$tableShape.Table.ApplyStyle($targetStyleId)
